$d = $word.ActiveDocument
$cursor = 0

function Find-From([int]$start, [string]$searchText) {
    $r = $d.Range($start, $d.Content.End)
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "NOT FOUND: $searchText"
        return $null
    }
    return $r
}

# [1] replace
$r = Find-From $cursor "• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations"
if ($r -ne $null) {
    $r.Text = "• Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions"
    $cursor = $r.End
}

# [2] replace
$r = Find-From $cursor "• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics"
if ($r -ne $null) {
    $r.Text = "• Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis"
    $cursor = $r.End
}

# [3] replace
$r = Find-From $cursor "• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets"
if ($r -ne $null) {
    $r.Text = "• Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets"
    $cursor = $r.End
}

# [4] replace
$r = Find-From $cursor "• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering"
if ($r -ne $null) {
    $r.Text = "• Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering"
    $cursor = $r.End
}

# [5] replace
$r = Find-From $cursor "• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications"
if ($r -ne $null) {
    $r.Text = "• Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications"
    $cursor = $r.End
}

# [6] replace
$r = Find-From $cursor "• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices"
if ($r -ne $null) {
    $r.Text = "• Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices"
    $cursor = $r.End
}

# [7] replace
$r = Find-From $cursor "• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES"
if ($r -ne $null) {
    $r.Text = "• Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES"
    $cursor = $r.End
}

# [8] replace
$r = Find-From $cursor "• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions"
if ($r -ne $null) {
    $r.Text = "• Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions"
    $cursor = $r.End
}

# [9] replace
$r = Find-From $cursor "• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI"
if ($r -ne $null) {
    $r.Text = "• Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI"
    $cursor = $r.End
}

# [10] replace
$r = Find-From $cursor "• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company`'s distinguishing products"
if ($r -ne $null) {
    $r.Text = "• Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company`'s distinguishing products"
    $cursor = $r.End
}

# [11] replace
$r = Find-From $cursor "• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices"
if ($r -ne $null) {
    $r.Text = "• Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices"
    $cursor = $r.End
}

# [12] replace
$r = Find-From $cursor "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research"
if ($r -ne $null) {
    $r.Text = "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions"
    $cursor = $r.End
}

# [13] replace
$r = Find-From $cursor "• Managed critical research operations for political campaigns"
if ($r -ne $null) {
    $r.Text = "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls"
    $cursor = $r.End
}

# [14] replace
$r = Find-From $cursor "• Conducted comprehensive polling and demographic analysis"
if ($r -ne $null) {
    $r.Text = "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"
    $cursor = $r.End
}

# [15] replace
$r = Find-From $cursor "• Developed strategic recommendations based on data analysis"
if ($r -ne $null) {
    $r.Text = "• Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
    $cursor = $r.End
}

# [16] replace
$r = Find-From $cursor "• Led research team in support of progressive political initiatives"
if ($r -ne $null) {
    $r.Text = "• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs"
    $cursor = $r.End
}

# [17] insert_after (x1) following cursor position
if ($r -ne $null) {
    $insPoint = $r
    $insPoint.Collapse(0)
    $insPoint.InsertParagraphAfter()
    $newPara = $d.Range($insPoint.End + 1, $insPoint.End + 1)
    $newPara.Text = "• Managed comprehensive research operations for progressive political initiatives and candidates"
    $insPoint = $newPara
    $cursor = $insPoint.End
    $r = $insPoint
}

# [18] replace
$r = Find-From $cursor "Political Research and Data Analysis"
if ($r -ne $null) {
    $r.Text = "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"
    $cursor = $r.End
}

# [19] replace
$r = Find-From $cursor "• Developed data analysis tools for political polling and research"
if ($r -ne $null) {
    $r.Text = "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
    $cursor = $r.End
}

# [20] replace
$r = Find-From $cursor "• Built statistical models for voter behavior analysis"
if ($r -ne $null) {
    $r.Text = "• Developed system that later became the Polling Consortium Database at The Analyst Institute"
    $cursor = $r.End
}

# [21] replace
$r = Find-From $cursor "• Created data visualization tools for research presentations"
if ($r -ne $null) {
    $r.Text = "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions"
    $cursor = $r.End
}

# [22] replace
$r = Find-From $cursor "• Supported senior researchers with technical analysis and reporting"
if ($r -ne $null) {
    $r.Text = "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle"
    $cursor = $r.End
}

# [23] insert_after (x2) following cursor position
if ($r -ne $null) {
    $insPoint = $r
    $insPoint.Collapse(0)
    $insPoint.InsertParagraphAfter()
    $newPara = $d.Range($insPoint.End + 1, $insPoint.End + 1)
    $newPara.Text = "• Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps"
    $insPoint = $newPara
    $insPoint.Collapse(0)
    $insPoint.InsertParagraphAfter()
    $newPara = $d.Range($insPoint.End + 1, $insPoint.End + 1)
    $newPara.Text = "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
    $insPoint = $newPara
    $cursor = $insPoint.End
    $r = $insPoint
}

# [24] replace
$r = Find-From $cursor "Political Field Operations and Data Management"
if ($r -ne $null) {
    $r.Text = "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"
    $cursor = $r.End
}

# [25] replace
$r = Find-From $cursor "• Managed field operations for political campaigns and research projects"
if ($r -ne $null) {
    $r.Text = "• Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions"
    $cursor = $r.End
}

# [26] replace
$r = Find-From $cursor "• Developed data collection and management systems for field work"
if ($r -ne $null) {
    $r.Text = "• Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm"
    $cursor = $r.End
}

# [27] replace
$r = Find-From $cursor "• Trained field staff on data collection protocols and quality control"
if ($r -ne $null) {
    $r.Text = "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
    $cursor = $r.End
}

# [28] replace
$r = Find-From $cursor "• Analyzed field data to inform campaign strategy and research findings"
if ($r -ne $null) {
    $r.Text = "• Created custom reports and data visualizations based on specific client requirements"
    $cursor = $r.End
}

# [29] insert_after (x2) following cursor position
if ($r -ne $null) {
    $insPoint = $r
    $insPoint.Collapse(0)
    $insPoint.InsertParagraphAfter()
    $newPara = $d.Range($insPoint.End + 1, $insPoint.End + 1)
    $newPara.Text = "• Introduced mapping and geospatial analysis into standard reporting procedures"
    $insPoint = $newPara
    $insPoint.Collapse(0)
    $insPoint.InsertParagraphAfter()
    $newPara = $d.Range($insPoint.End + 1, $insPoint.End + 1)
    $newPara.Text = "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
    $insPoint = $newPara
    $cursor = $insPoint.End
    $r = $insPoint
}

